$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": a new handoff xliff was generated for
# eceeb463-15c4-474c-b45c-71333f559960.md (row 7 of each table), updating the
# "Latest Handoff Datetime" for zh-cn and de-de, and the overall
# "Latest HO Xliff Generate Date" on the Overview sheet.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G7 -> "Latest HO Xliff Generate Date" for eceeb463...md
$wsOverview.Range("G7").Value = "2016-08-21 10:48:08"

# zh-cn!H7 -> "Latest Handoff Datetime" for eceeb463...zh-cn.xlf
$wsZhCn.Range("H7").Value = "2016-08-21 10:48:00"

# de-de!H7 -> "Latest Handoff Datetime" for eceeb463...de-de.xlf
$wsDeDe.Range("H7").Value = "2016-08-21 10:48:08"
